# Updates the cryptos list (Price / Volume(1h) columns, plus a couple of
# coin-name/link/price swaps) to match the latest scrape.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> @{ Column = NewValue } for cells that changed.
# Columns: B=Coin, C=Link, D=Price, E=Volume(1h)
$changes = @{
    2  = @{ D = "23.195.23";  E = "  +0.30%  " }
    3  = @{ D = "1.602.44";   E = "  -0.14%  " }
    4  = @{ D = "1.000";      E = "  -0.05%  " }
    5  = @{ D = "1.000";      E = "  -0.01%  " }
    6  = @{ D = "303.21";     E = "  +0.37%  " }
    7  = @{ D = "0.3783";     E = "  -0.08%  " }
    8  = @{ D = "51.95";      E = "  +3.34%  " }
    9  = @{ D = "0.3623";     E = "  -1.19%  " }
    10 = @{ D = "1.270";      E = "  -0.41%  " }
    11 = @{ D = "1.000";      E = "  -0.04%  " }
    12 = @{ D = "0.08125";    E = "  -0.33%  " }
    13 = @{ D = "22.83";      E = "  -0.16%  " }
    14 = @{ D = "6.603";      E = "  -0.43%  " }
    15 = @{ E = "  +0.01%  " }
    16 = @{ D = "0.00001245"; E = "  -1.34%  " }
    17 = @{ D = "1.604.40";   E = "  +0.01%  " }
    18 = @{ D = "93.98";      E = "  +1.92%  " }
    19 = @{ D = "0.06882";    E = "  +0.13%  " }
    20 = @{ D = "18.08";      E = "  -1.36%  " }
    21 = @{ D = "6.542";      E = "  -0.96%  " }
    22 = @{ E = "  +0.01%  " }
    23 = @{ D = "12.98";      E = "  -0.97%  " }
    24 = @{ D = "23.204.30";  E = "  +0.31%  " }
    25 = @{ B = "LidoDAOToken"; C = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"; D = "3.020"; E = "  +7.50%  " }
    26 = @{ B = "Toncoin";      C = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton";      D = "2.393"; E = "  +1.47%  " }
    27 = @{ E = "  +0.17%  " }
    28 = @{ D = "149.89";     E = "  -0.34%  " }
    29 = @{ D = "5.236";      E = "  -0.59%  " }
    30 = @{ D = "133.81";     E = "  -0.18%  " }
    31 = @{ D = "2.360";      E = "  -0.84%  " }
    32 = @{ D = "6.753";      E = "  -1.80%  " }
    33 = @{ D = "1.782.14";   E = "  +0.34%  " }
    34 = @{ D = "0.9641";     E = "  +0.10%  " }
    35 = @{ D = "0.07494";    E = "  -3.20%  " }
    36 = @{ B = "VeChain";   C = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"; D = "0.02724"; E = "  -0.31%  " }
    37 = @{ B = "FraxShare"; C = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs";    D = "10.24";   E = "  -2.59%  " }
    38 = @{ D = "0.2524";     E = "  -1.28%  " }
    39 = @{ D = "0.08796";    E = "  -1.24%  " }
    40 = @{ D = "6.089";      E = "  -3.38%  " }
    41 = @{ D = "0.7113";     E = "  +0.16%  " }
    42 = @{ D = "1.364";      E = "  -0.50%  " }
    43 = @{ D = "12.51";      E = "  -1.24%  " }
    44 = @{ D = "15.72";      E = "  +2.09%  " }
    45 = @{ D = "0.6549";     E = "  -1.50%  " }
    46 = @{ D = "2.315";      E = "  -0.45%  " }
    47 = @{ D = "4.017";      E = "  +0.31%  " }
    48 = @{ D = "132.25";     E = "  -0.36%  " }
    49 = @{ D = "0.07950";    E = "  +0.01%  " }
    50 = @{ D = "1.206";      E = "  -3.93%  " }
    51 = @{ E = "  -0.34%  " }
}

foreach ($row in $changes.Keys) {
    $cellChanges = $changes[$row]
    foreach ($col in $cellChanges.Keys) {
        $address = "$col$row"
        $range = $ws.Range($address)
        $value = $cellChanges[$col]

        if ($col -eq "D") {
            # Price values look numeric (e.g. "1.000", "303.21") but must be
            # stored as plain text, exactly as scraped. Force a text format
            # before assignment so Excel doesn't silently convert/round them
            # to a real number, then restore the default "Normal" style so
            # no stray number-format style lingers on the cell.
            $range.NumberFormat = "@"
            $range.Value = $value
            $range.Style = "Normal"
        } else {
            $range.Value = $value
        }
    }
}
